$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 63
$values = @(
    "01-04-2021",
    40018334,
    31863485,
    31465498,
    397988,
    0,
    0,
    1302,
    0,
    2,
    8153547,
    46340873,
    44954,
    0,
    0,
    0,
    13286836,
    -6322539,
    0
)

$ws.Cells.Item($row, 1).Value = "'" + $values[0]
$ws.Cells.Item($row, 1).ClearFormats()
for ($i = 1; $i -lt $values.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item($row, $col).Value = $values[$i]
}
